# ECU-06.docx edit:
#  1. Remove the <w:bookmarkStart/bookmarkEnd w:name="_GoBack"/> pair that
#     sits right after "Descripcion general".
#  2. Split the single run that holds the "Pre-condiciones" body text into
#     five runs (same Arial rPr on every run) and re-insert a fresh
#     _GoBack bookmark pair at the end of that paragraph.

$d = $word.ActiveDocument

# --- Step 1: drop the old _GoBack bookmark -------------------------------
$bookmarks = $d.Bookmarks
if ($bookmarks.Exists("_GoBack")) {
    $bookmarks.Item("_GoBack").Delete()
}

# --- Step 2: locate the "Pre-condiciones" body paragraph -----------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -like "*funcionando correctamente y el MVZ que se va a registrar no debe existir en el sistema*") {
        $targetPara = $candidate
        break
    }
}

if ($targetPara -eq $null) {
    throw "Could not find the Pre-condiciones paragraph"
}

$rng = $targetPara.Range

# Build the replacement paragraph: same pPr/paragraph marks, but the body
# text is now spread across five runs (identical Arial formatting) and a
# new _GoBack bookmark wraps the very end of the paragraph.
$newParaXml = '<w:p w14:paraId="3A5A5BFF" w14:textId="098407C1" w:rsidR="00A77B3E" w:rsidRPr="004D7E7A" w:rsidRDefault="00A77B3E" w:rsidP="00245D75"><w:pPr><w:pStyle w:val="Textodecuerpo"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>E</w:t></w:r><w:r w:rsidR="005433B3"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>l sistema debe estar funcionando correctamente</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>, el MVZ debe estar autentificado correctamente</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> y el </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">nuevo </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>MVZ que se va a registrar no debe existir en el sistema.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$xmlPackage = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + '<w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xmlPackage, $null)
